$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.268.88"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "2.249.16"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "2.591.15"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.837"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.236.51"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "44.077.31"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "0.0₃0970"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.88%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.107"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  -6.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.95%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.23%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "1.760.61"
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.62%  "
